# Nexial "#system" command-reference sheet update:
#  - add new `tn.5250` command category (close/open/saveText/typeKeys/updateScreenFields)
#  - add new `ocr(image,saveVar)` command to the existing `image` category
#  - rename `colorbit(source,bit,saveTo)` -> `colorbit(image,bit,saveTo)` for consistency
#  - keep all the defined-name ranges in sync with the shifted columns/rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("#system")

# ---------------------------------------------------------------------------
# 1) New category column: insert a whole column before Z (this correctly
#    shifts Z:AE -> AA:AF for every row, column-insert only affects this
#    sheet) and populate the new Z column with the "tn.5250" category plus
#    its five commands.
# ---------------------------------------------------------------------------
$ws.Columns("Z").Insert()

$ws.Range("Z1").Value = "tn.5250"
$ws.Range("Z2").Value = "close(profile)"
$ws.Range("Z3").Value = "open(profile)"
$ws.Range("Z4").Value = "saveText(profile,var)"
$ws.Range("Z5").Value = "typeKeys(profile,keystrokes)"
$ws.Range("Z6").Value = "updateScreenFields(profile)"

# ---------------------------------------------------------------------------
# 2) Register the new category name in column A ("target" list), which is
#    sorted alphabetically: "tn.5250" sits between "step" and "web".
#    NOTE: Range.Insert() on a single cell in this sheet shifts the whole
#    row, not just the column, so the shift is done manually cell-by-cell
#    (bottom-up, to avoid clobbering values) instead of via Insert().
# ---------------------------------------------------------------------------
for ($r = 31; $r -ge 26; $r--) {
    $ws.Range("A" + ($r + 1)).Value = $ws.Range("A$r").Value()
}
$ws.Range("A26").Value = "tn.5250"

# ---------------------------------------------------------------------------
# 3) Update the "image" command category: rename the `colorbit` signature and
#    insert the new `ocr(image,saveVar)` command (alphabetically, right
#    before `resize`). Same manual bottom-up shift, restricted to column K.
# ---------------------------------------------------------------------------
$ws.Range("K2").Value = "colorbit(image,bit,saveTo)"

$ws.Range("K8").Value = $ws.Range("K7").Value()
$ws.Range("K7").Value = $ws.Range("K6").Value()
$ws.Range("K6").Value = "ocr(image,saveVar)"

# ---------------------------------------------------------------------------
# 4) Keep the workbook-level defined names in sync with the shifted ranges.
# ---------------------------------------------------------------------------
$wb.Names("image").RefersTo     = "='#system'!`$K`$2:`$K`$8"
$wb.Names("target").RefersTo    = "='#system'!`$A`$2:`$A`$32"
$wb.Names("web").RefersTo       = "='#system'!`$AA`$2:`$AA`$144"
$wb.Names("webalert").RefersTo  = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names("webcookie").RefersTo = "='#system'!`$AC`$2:`$AC`$10"
$wb.Names("ws").RefersTo        = "='#system'!`$AD`$2:`$AD`$17"
$wb.Names("ws.async").RefersTo  = "='#system'!`$AE`$2:`$AE`$8"
$wb.Names("xml").RefersTo       = "='#system'!`$AF`$2:`$AF`$27"
$wb.Names.Add("tn.5250", "='#system'!`$Z`$2:`$Z`$6")
